$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on target cells so numeric-looking strings (e.g. "143.01",
# "57.922.37", percentages) are preserved verbatim instead of being coerced to
# floating point numbers by the COM value setter.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.922.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.573.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.88%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.25"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.01"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.39%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.591.48"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.62"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.324"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.77%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.031.38"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.893.59"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.30"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.93%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.31%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.579.55"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "338.85"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.21"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.32"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.46"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.402"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.684.50"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.97"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0750"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.29"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.70%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.65"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.40"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.01"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.14"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.872"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.18"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.833"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.53"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "270.81"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.66"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0954"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.588"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.77"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0522"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.976.07"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.39"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.20%  "
